$d = $word.ActiveDocument

# --- Change 1 -----------------------------------------------------------------
# " assessment author, in consultation with the" -> " assessment developer, in
# consultation with the" (i.e. "author" becomes "developer").
# "author," (with the trailing comma) is unique in the document - the bare word
# "author" also shows up in "assessment author(s)" elsewhere in the checklist,
# so searching on the comma avoids touching that unrelated occurrence. Replace
# type 1 = wdReplaceOne, just to be extra safe even though the text is unique.
$null = $d.Content.Find.Execute("author,", $false, $false, $false, $false, $false, $true, 1, $false, "developer,", 1)

# --- Change 2 -----------------------------------------------------------------
# The hyperlinked text "House Style for Site Text " + "a" + "nd Zotero Entries"
# (three runs, first occurrence only - inside the "house style as set out in"
# sentence) becomes a single run "House Style for Site Text and Zotero Entries".
$rng = $d.Content.Duplicate
$rng.Find.ClearFormatting()
$rng.Find.Text = "House Style for Site Text and Zotero Entries"
$null = $rng.Find.Execute()
if ($rng.Find.Found) {
    $start = $rng.Start
    $end = $rng.End
    # Setting identical text is a no-op in this engine, so briefly swap in a
    # placeholder to force the real text mutation (which collapses the three
    # runs found above into one), then write the final text back.
    $tmp = $d.Range($start, $end)
    $tmp.Text = "~~~"
    $tmp2 = $d.Range($start, $start + 3)
    $tmp2.Text = "House Style for Site Text and Zotero Entries"
    # The merge picks up the plain paragraph formatting rather than the
    # hyperlink's own blue/underline run formatting, so restore it explicitly.
    $newLen = "House Style for Site Text and Zotero Entries".Length
    $fixed = $d.Range($start, $start + $newLen)
    $fixed.Font.Color = 0xC16305
    $fixed.Font.Underline = 1
}

# --- Change 3 -----------------------------------------------------------------
# "[" -> "(" and "]" -> ")" around "an assessment and corresponding essay".
# Both characters are unique in the document.
$null = $d.Content.Find.Execute("[", $false, $false, $false, $false, $false, $true, 1, $false, "(", 2)
$null = $d.Content.Find.Execute("]", $false, $false, $false, $false, $false, $true, 1, $false, ")", 2)

Write-Output "done"
